$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the edited range keeps a text format so numeric-looking
# strings (e.g. "0", "1") are stored as text rather than numbers,
# matching the original data's string typing.
$ws.Range("A2:C3").NumberFormat = "@"

# Update row 2
$ws.Range("A2").Value = "HSA Advocates"
$ws.Range("B2").Value = "17sec"
$ws.Range("C2").Value = "0"

# Update row 3
$ws.Range("A3").Value = "Cobalt Legal"
$ws.Range("B3").Value = "31sec"
$ws.Range("C3").Value = "1"

# Clear out rows 4 through 52 (columns A:C)
$ws.Range("A4:C52").ClearContents()
